$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The document contains two near-identical "CUU" blocks (one in the Italian
# section, one in the English section), each with the text:
#   "Codice univoco ufficio (CUU): 4S488Q."
# We need to wrap the "4S488Q" run in each block with a new bookmark:
#   - bookmark_cuu      (Italian / first occurrence)
#   - bookmark_cuu_en   (English / second occurrence)
# ---------------------------------------------------------------------------

# --- First occurrence (Italian section) ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("4S488Q", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found1) {
    $d.Bookmarks.Add("bookmark_cuu", $r1)
}

# --- Second occurrence (English section) ---
# Search the remainder of the document, starting right after the first match,
# so we land on the second "4S488Q" occurrence.
$r2 = $d.Range($r1.End, $d.Content.End)
$found2 = $r2.Find.Execute("4S488Q", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found2) {
    $d.Bookmarks.Add("bookmark_cuu_en", $r2)
}
